# Helper: write a value to a Range while forcing Excel to keep it as TEXT
# (prevents "33.68" / "001481" style values from being auto-coerced to numbers),
# without leaving a residual NumberFormat/style change behind on the cell.
function Set-TextValue {
    param($range, $value)
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计" (so the tab
#    order becomes 总计, 2022-Q4, 2021-Q3, 2021-Q2).
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# ------------------------------------------------------------------
# 2. Populate the new "2022-Q4" sheet with the fund holdings table.
# ------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
Set-TextValue $q4.Range("B2") "001481"
$q4.Range("C2").Value = "华宝油气（QDII）美元"
Set-TextValue $q4.Range("D2") "50.10"
Set-TextValue $q4.Range("E2") "94.65"
Set-TextValue $q4.Range("F2") "2.11"
Set-TextValue $q4.Range("G2") "1.0571"
$q4.Range("H2").Value = 8

$q4.Range("A3").Value = 1
Set-TextValue $q4.Range("B3") "007844"
$q4.Range("C3").Value = "华宝油气（QDII）人民币 C"
Set-TextValue $q4.Range("D3") "27.91"
Set-TextValue $q4.Range("E3") "94.65"
Set-TextValue $q4.Range("F3") "2.11"
Set-TextValue $q4.Range("G3") "0.5889"
$q4.Range("H3").Value = 8

$q4.Range("A4").Value = 2
Set-TextValue $q4.Range("B4") "162411"
$q4.Range("C4").Value = "华宝油气（QDII）人民币A"
Set-TextValue $q4.Range("D4") "22.19"
Set-TextValue $q4.Range("E4") "94.65"
Set-TextValue $q4.Range("F4") "2.11"
Set-TextValue $q4.Range("G4") "0.4682"
$q4.Range("H4").Value = 8

# ------------------------------------------------------------------
# 3. Update the "总计" summary sheet: a new row for 2022-Q4 is
#    inserted at the top of the data, pushing 2021-Q3 / 2021-Q2 down.
# ------------------------------------------------------------------
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("D2").Value = 2.11

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q3"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 2.09

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q2"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 2.17
